$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Table 0")

$data = @(
    @(2, "193,83 ", "199,46 ", "133,28 "),
    @(3, "144,93 ", "147,33 ", "80,85 "),
    @(4, "155,13 ", "156,93 ", "77,87 "),
    @(5, "225,95 ", "234,93 ", "118,90 "),
    @(6, "172,15 ", "167,93 ", "78,58 "),
    @(7, "211,05 ", "198,05 ", "126,84 "),
    @(8, "194,54 ", "184,82 ", "125,44 "),
    @(9, "190,55 ", "182,00 ", "113,26 "),
    @(10, "213,28 ", "204,37 ", "114,31 "),
    @(11, "243,31 ", "216,46 ", "131,74 "),
    @(12, "212,81 ", "184,35 ", "113,14 "),
    @(13, "178,02 ", "167,48 ", "98,38 "),
    @(14, "177,31 ", "178,39 ", "105,52 "),
    @(15, "205,66 ", "199,34 ", "112,44 "),
    @(16, "223,82 ", "212,81 ", "85,85 "),
    @(17, "173,92 ", "187,04 ", "91,70 "),
    @(18, "189,97 ", "187,04 ", "106,34 "),
    @(19, "151,28 ", "139,86 ", "87,53 "),
    @(20, "243,26 ", "217,84 ", "104,59 "),
    @(21, "169,64 ", "171,84 ", "111,07 "),
    @(22, "211,87 ", "208,71 ", "110,44 "),
    @(23, "173,25 ", "176,79 ", "82,15 "),
    @(24, "183,98 ", "200,03 ", "106,33 "),
    @(25, "180,95 ", "187,39 ", "107,63 "),
    @(26, "185,17 ", "181,54 ", "90,65 "),
    @(27, "148,48 ", "155,09 ", "83,56 ")
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
}

$ws.Range("F4").Select()